$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '72.041.49'
$ws.Range("E2").Value2 = '  +0.37%  '
$ws.Range("D3").Value2 = '4.043.40'
$ws.Range("E3").Value2 = '  -0.11%  '
$ws.Range("E4").Value2 = '  -0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = '537.62'
$c.Style = "Normal"
$ws.Range("E5").Value2 = '  +0.57%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = '149.49'
$c.Style = "Normal"
$ws.Range("E6").Value2 = '  -2.32%  '
$ws.Range("D7").Value2 = '4.037.42'
$ws.Range("E7").Value2 = '  -0.06%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = '0.699'
$c.Style = "Normal"
$ws.Range("E8").Value2 = '  +0.68%  '
$ws.Range("E9").Value2 = '  -0.04%  '
$ws.Range("E10").Value2 = '  -1.19%  '
$ws.Range("E11").Value2 = '  -2.36%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = '53.39'
$c.Style = "Normal"
$ws.Range("E12").Value2 = '  +9.20%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = '0.0000330'
$c.Style = "Normal"
$ws.Range("E13").Value2 = '  -0.71%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = '10.93'
$c.Style = "Normal"
$ws.Range("E14").Value2 = '  -0.30%  '
$ws.Range("D15").Value2 = '4.680.95'
$ws.Range("E15").Value2 = '  -0.29%  '
$ws.Range("D16").Value2 = '4.035.71'
$ws.Range("E16").Value2 = '  -0.16%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = '14.31'
$c.Style = "Normal"
$ws.Range("E17").Value2 = '  -0.87%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = '20.72'
$c.Style = "Normal"
$ws.Range("E18").Value2 = '  -1.54%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = '1.21'
$c.Style = "Normal"
$ws.Range("E19").Value2 = '  -1.27%  '
$ws.Range("E20").Value2 = '  -1.21%  '
$ws.Range("D21").Value2 = '72.035.41'
$ws.Range("E21").Value2 = '  +0.34%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = '438.01'
$c.Style = "Normal"
$ws.Range("E22").Value2 = '  +0.58%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = '98.17'
$c.Style = "Normal"
$ws.Range("E23").Value2 = '  -1.34%  '
$ws.Range("E24").Value2 = '  -5.60%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = '4.27'
$c.Style = "Normal"
$ws.Range("E25").Value2 = '  +1.11%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = '14.56'
$c.Style = "Normal"
$ws.Range("E26").Value2 = '  -1.75%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = '4.33'
$c.Style = "Normal"
$ws.Range("E27").Value2 = '  +21.80%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = '11.33'
$c.Style = "Normal"
$ws.Range("E28").Value2 = '  -0.70%  '
$ws.Range("E29").Value2 = '  -1.98%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = '37.11'
$c.Style = "Normal"
$ws.Range("E31").Value2 = '  -0.46%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = '8.46'
$c.Style = "Normal"
$ws.Range("E32").Value2 = '  +24.60%  '
$ws.Range("E33").Value2 = '  +1.91%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = '50.20'
$c.Style = "Normal"
$ws.Range("E34").Value2 = '  +16.62%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = '13.59'
$c.Style = "Normal"
$ws.Range("E35").Value2 = '  -0.84%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = '677.63'
$c.Style = "Normal"
$ws.Range("E36").Value2 = '  -0.41%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value2 = '66.69'
$c.Style = "Normal"
$ws.Range("E37").Value2 = '  +0.09%  '
$ws.Range("E38").Value2 = '  +5.92%  '
$ws.Range("D39").Value2 = '0.0₃0876'
$ws.Range("E39").Value2 = '  +1.91%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = '3.45'
$c.Style = "Normal"
$ws.Range("E40").Value2 = '  +8.58%  '
$ws.Range("E41").Value2 = '  -6.24%  '
$ws.Range("E42").Value2 = '  -1.34%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = '11.26'
$c.Style = "Normal"
$ws.Range("E43").Value2 = '  +17.57%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = '0.998'
$c.Style = "Normal"
$ws.Range("E44").Value2 = '  -0.15%  '
$ws.Range("B45").Value2 = 'VeChain'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = '0.0494'
$c.Style = "Normal"
$ws.Range("E45").Value2 = '  -0.94%  '
$ws.Range("B46").Value2 = 'FirstDigitalUSD'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = '0.999'
$c.Style = "Normal"
$ws.Range("E46").Value2 = '  +0.05%  '
$ws.Range("E47").Value2 = '  -0.88%  '
$ws.Range("E48").Value2 = '  -2.14%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = '3.12'
$c.Style = "Normal"
$ws.Range("E49").Value2 = '  +1.84%  '
$ws.Range("E50").Value2 = '  -3.29%  '
$ws.Range("D51").Value2 = '2.843.69'
$ws.Range("E51").Value2 = '  +8.70%  '
